# Re-delete the 'culture_collection' column (column R) from the MIGS
# BioSample template, matching INSDC2017 review (column was previously
# re-added, this removes it again).
#
# Deleting the column shifts the worksheet cell data and the shared
# strings table automatically, but cell comments (and their underlying
# VML note shapes) are anchored to fixed cell references and are not
# re-targeted by a column delete. So the comment text attached to every
# column from R (culture_collection) onward is first shifted left by
# one position (absorbing the text that belonged to the next column to
# the right), the now-duplicated trailing comment is removed, and only
# then is the now-empty culture_collection column itself deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns R..AG will each take on the comment text that currently
# belongs to the next column over (S..AH), i.e. the comment that used
# to describe 'culture_collection' (R15) becomes the 'encoded_traits'
# description (previously S15), and so on through the last column.
$destCols = @("R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG")
$srcCols  = @("S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH")

for ($i = 0; $i -lt $destCols.Length; $i++) {
    $srcCell = $srcCols[$i] + "15"
    $destCell = $destCols[$i] + "15"
    $text = $ws.Range($srcCell).Comment.Text()
    [void]$ws.Range($destCell).Comment.Text($text)
}

# The trailing column's (AH15) comment text has now been copied onto
# AG15, so the original AH15 comment is a leftover duplicate - remove it.
[void]$ws.Range("AH15").Comment.Delete()

# Finally, delete the (now textually-vacated) culture_collection column
# itself; this shifts all worksheet values/shared-string references one
# column to the left, matching the comment shift performed above.
[void]$ws.Range("R15").EntireColumn.Delete()
